# Domino JTAG SPI BOM: Rev. B -> Rev. C

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Domino JTAG SPI Rev. B")

# Rename the sheet (this also updates the Print_Area_* defined names that
# reference the sheet by name).
$ws.Name = "Domino JTAG SPI Rev. C"

# The plain "Print_Area" defined name isn't auto-retargeted by the rename,
# so fix it up explicitly.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Domino JTAG SPI Rev. C!Print_Area") {
        $n.RefersTo = "='Domino JTAG SPI Rev. C'!`$A`$1:`$I`$12"
    }
}

# Row 8: R1 changes from a 1k resistor to a 0R jumper.
$ws.Range("E8").Value = "R0402_0R_5%_62.5mW"
$ws.Range("H8").Value = "RES 0.0 OHM 1/16W JUMP 0402 SMD"

# Row 9: add R16 to the 10k group, bumping the quantity to 7.
$ws.Range("B9").Value = 7
$ws.Range("G9").Value = "R2, R3, R4, R5, R7, R8, R16"

# Selection moved to G10 in the saved file.
$ws.Range("G10").Select()
